# Reorder the player roster rows (A2:C19) to the new order while keeping
# the same header row and the same set of (Player, Position, Team) tuples.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Dejounte Murray",    "PG,SG",   "New Orleans Pelicans"),
    @("Russell Westbrook",  "PG,SG",   "Denver Nuggets"),
    @("Tyus Jones",         "PG",      "Phoenix Suns"),
    @("Jalen Green",        "PG,SG",   "Houston Rockets"),
    @("Jaylen Brown",       "SG,SF",   "Boston Celtics"),
    @("Draymond Green",     "PF,C",    "Golden State Warriors"),
    @("Khris Middleton",    "SF",      "Milwaukee Bucks"),
    @("Deni Avdija",        "SF,PF",   "Portland Trail Blazers"),
    @("Rudy Gobert",        "C",       "Minnesota Timberwolves"),
    @("Nikola Jokic",       "C",       "Denver Nuggets"),
    @("Jakob Poeltl",       "C",       "Toronto Raptors"),
    @("Jonas Valanciunas",  "C",       "Washington Wizards"),
    @("Keon Johnson",       "PG,SG",   "Brooklyn Nets"),
    @("Chris Paul",         "PG",      "San Antonio Spurs"),
    @("Pascal Siakam",      "SF,PF,C", "Indiana Pacers"),
    @("Paolo Banchero",     "SF,PF",   "Orlando Magic"),
    @("Chet Holmgren",      "PF,C",    "Oklahoma City Thunder"),
    @("Jalen Suggs",        "PG,SG",   "Orlando Magic")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
